$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '246.62'
$ws.Range('D2').Style = 'Normal'

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.83'
$ws.Range('D3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.471'
$ws.Range('D4').Style = 'Normal'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.05742'
$ws.Range('D5').Style = 'Normal'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '3.435'
$ws.Range('D6').Style = 'Normal'

$ws.Range('B7').NumberFormat = '@'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('B7').Style = 'Normal'

$ws.Range('C7').NumberFormat = '@'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('C7').Style = 'Normal'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '6.324'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '6KuCoinTokenKCS'
$ws.Range('E7').Style = 'Normal'

$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('B8').Style = 'Normal'

$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C8').Style = 'Normal'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.8114'
$ws.Range('D8').Style = 'Normal'

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('E8').Style = 'Normal'

$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('B9').Style = 'Normal'

$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('C9').Style = 'Normal'

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8834'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('E9').Style = 'Normal'

$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'One'
$ws.Range('B10').Style = 'Normal'

$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('C10').Style = 'Normal'

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.01095'
$ws.Range('D10').Style = 'Normal'

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '9OneONEBestin24h'
$ws.Range('E10').Style = 'Normal'

$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('B11').Style = 'Normal'

$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('C11').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1443'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('E11').Style = 'Normal'

$ws.Range('B12').NumberFormat = '@'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('B12').Style = 'Normal'

$ws.Range('C12').NumberFormat = '@'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('C12').Style = 'Normal'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07394'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('E12').Style = 'Normal'

$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('B13').Style = 'Normal'

$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('C13').Style = 'Normal'

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03009'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('E13').Style = 'Normal'

$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('B14').Style = 'Normal'

$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('C14').Style = 'Normal'

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03077'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('E14').Style = 'Normal'

$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('B15').Style = 'Normal'

$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('C15').Style = 'Normal'

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09405'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('E15').Style = 'Normal'

$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('B16').Style = 'Normal'

$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('C16').Style = 'Normal'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.924'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('E16').Style = 'Normal'

$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('B17').Style = 'Normal'

$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('C17').Style = 'Normal'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.001588'
$ws.Range('D17').Style = 'Normal'

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('E17').Style = 'Normal'

$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('B18').Style = 'Normal'

$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('C18').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04802'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('E18').Style = 'Normal'

$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('B19').Style = 'Normal'

$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('C19').Style = 'Normal'

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006218'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('E19').Style = 'Normal'

$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('B20').Style = 'Normal'

$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('C20').Style = 'Normal'

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.005123'
$ws.Range('D20').Style = 'Normal'

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('E20').Style = 'Normal'

$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('B21').Style = 'Normal'

$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('C21').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0009968'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '20BitKanKAN'
$ws.Range('E21').Style = 'Normal'

$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('B22').Style = 'Normal'

$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('C22').Style = 'Normal'

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001499'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('E22').Style = 'Normal'

$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'LEO'
$ws.Range('B23').Style = 'Normal'

$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C23').Style = 'Normal'

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.750'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '22LEOLEO'
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.190'
$ws.Range('D24').Style = 'Normal'

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3280'
$ws.Range('D25').Style = 'Normal'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0002985'
$ws.Range('D27').Style = 'Normal'

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03905'
$ws.Range('D40').Style = 'Normal'

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006784'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1071'
$ws.Range('D42').Style = 'Normal'

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002898'
$ws.Range('D43').Style = 'Normal'

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.007370'
$ws.Range('D44').Style = 'Normal'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005634'
$ws.Range('D45').Style = 'Normal'

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000749'
$ws.Range('D46').Style = 'Normal'

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.3797'
$ws.Range('D47').Style = 'Normal'

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1705'
$ws.Range('D48').Style = 'Normal'

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00002098'
$ws.Range('D49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.01009'
$ws.Range('D50').Style = 'Normal'
